# The workbook originally has a sheet named "strategy_id-5008".
# This edit renames it to "strategy_id-5007", and adds a duplicate of it
# (same data/layout) right after it, named "strategy_id-5009".

$wb = $excel.ActiveWorkbook

# Rename the existing "strategy_id-5008" sheet to "strategy_id-5007".
$ws = $wb.Worksheets.Item("strategy_id-5008")
$ws.Name = "strategy_id-5007"

# Duplicate it, placing the copy immediately after the original,
# and name the new sheet "strategy_id-5009".
$ws.Copy($null, $ws)
$newWs = $wb.Worksheets.Item($ws.Index + 1)
$newWs.Name = "strategy_id-5009"
